# Insert a new weekly price record for Espinaca (Agrícola del Norte S.A. de Arica)
# This shifts the existing rows 69-96 down to 70-97 and populates the new row 69
# with the latest observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 69 (old row 69 and everything below moves down by one)
$ws.Rows.Item(69).Insert()

$ws.Cells.Item(69, 1).Value = 1
$ws.Cells.Item(69, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(69, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(69, 4).Value = 45006
$ws.Cells.Item(69, 5).Value = 15
$ws.Cells.Item(69, 6).Value = 100112012
$ws.Cells.Item(69, 7).Value = "Espinaca"
$ws.Cells.Item(69, 8).Value = "Sin especificar"
$ws.Cells.Item(69, 9).Value = "Primera"
$ws.Cells.Item(69, 10).Value = 200
$ws.Cells.Item(69, 11).Value = 3500
$ws.Cells.Item(69, 12).Value = 4000
$ws.Cells.Item(69, 13).Value = 3750
$ws.Cells.Item(69, 14).Value = "`$/atado 2,5 a 3 kilos"
$ws.Cells.Item(69, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(69, 16).Value = 1250
$ws.Cells.Item(69, 17).Value = 3
$ws.Cells.Item(69, 18).Value = "Hortaliza"
